$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-26 Friday", "2024-07-27 Saturday"),
    @("524÷2=262, 0", "267÷5=53, 2"),
    @("434÷7=62, 0", "389÷3=129, 2"),
    @("372÷8=46, 4", "144÷7=20, 4"),
    @("744÷8=93, 0", "290÷6=48, 2"),
    @("489÷3=163, 0", "473÷7=67, 4"),
    @("502÷7=71, 5", "980÷5=196, 0"),
    @("874÷9=97, 1", "250÷8=31, 2"),
    @("592÷4=148, 0", "920÷7=131, 3"),
    @("581÷9=64, 5", "832÷2=416, 0"),
    @("226÷4=56, 2", "588÷2=294, 0"),
    @("737÷3=245, 2", "189÷9=21, 0"),
    @("692÷5=138, 2", "226÷8=28, 2"),
    @("103÷2=51, 1", "285÷2=142, 1"),
    @("767÷2=383, 1", "973÷8=121, 5"),
    @("983÷6=163, 5", "234÷4=58, 2"),
    @("761÷6=126, 5", "378÷9=42, 0"),
    @("881÷6=146, 5", "262÷4=65, 2"),
    @("510÷8=63, 6", "960÷4=240, 0"),
    @("324÷6=54, 0", "853÷4=213, 1"),
    @("292÷4=73, 0", "277÷3=92, 1"),
    @("981÷8=122, 5", "545÷5=109, 0"),
    @("469÷5=93, 4", "514÷4=128, 2"),
    @("608÷6=101, 2", "912÷6=152, 0"),
    @("862÷6=143, 4", "331÷9=36, 7"),
    @("700÷8=87, 4", "351÷8=43, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
